$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh re-ordered/replaced the entries in rows 4-10
# (dates, volumes and prices for "Arándano (blue)"). Each block below
# rewrites one row with the values it now holds after the refresh.

# Row 4: apply data originally from row 9
$ws.Range("D4").Value = 44596
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 2500
$ws.Range("O4").Value = 2700
$ws.Range("P4").Value = 2600
$ws.Range("R4").Value = 'Provincia de Linares'
$ws.Range("S4").Value = 1300

# Row 5: apply data originally from row 10
$ws.Range("D5").Value = 44594
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 2500
$ws.Range("O5").Value = 2800
$ws.Range("P5").Value = 2650
$ws.Range("S5").Value = 1325

# Row 6: apply data originally from row 4
$ws.Range("D6").Value = 44539
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 3800
$ws.Range("O6").Value = 4000
$ws.Range("P6").Value = 3900
$ws.Range("Q6").Value = '$/bandeja 2 kilos'
$ws.Range("R6").Value = 'Región del Maule'
$ws.Range("S6").Value = 1950
$ws.Range("T6").Value = 2

# Row 8: apply data originally from row 5
$ws.Range("D8").Value = 44187
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 2800
$ws.Range("O8").Value = 3000
$ws.Range("P8").Value = 2900
$ws.Range("S8").Value = 1450

# Row 9: apply data originally from row 6
$ws.Range("D9").Value = 44187
$ws.Range("M9").Value = 65
$ws.Range("N9").Value = 1400
$ws.Range("O9").Value = 1500
$ws.Range("P9").Value = 1446
$ws.Range("Q9").Value = '$/envase 1 kilo'
$ws.Range("R9").Value = 'Provincia de Diguillín'
$ws.Range("S9").Value = 1446
$ws.Range("T9").Value = 1

# Row 10: apply data originally from row 8
$ws.Range("D10").Value = 44174
$ws.Range("M10").Value = 150
$ws.Range("N10").Value = 3700
$ws.Range("O10").Value = 3800
$ws.Range("P10").Value = 3747
$ws.Range("S10").Value = 1874
